$wb = $excel.ActiveWorkbook

# --- Users sheet: add a new user row (daniel34 / Daniel34! / 777) ---
$users = $wb.Worksheets.Item("Users")
$users.Cells.Item(15, 1).Value = "daniel34"
$users.Cells.Item(15, 2).Value = "Daniel34!"
$users.Cells.Item(15, 3).Value = 777

# --- Birds sheet: add two new bird records ---
$birds = $wb.Worksheets.Item("Birds")

# Exact calendar date (midnight) used for both new records.
$recordDate = Get-Date -Year 2023 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0

# Row 13: American Gouldian / Central America / Male, cage a3, user 777
$birds.Cells.Item(13, 1).Value = 12
$birds.Cells.Item(13, 2).Value = "American Gouldian"
$birds.Cells.Item(13, 3).Value = "Central America"
$birds.Cells.Item(13, 4).Value = "Male"
# Copy the date formatting from an existing date cell so the new cell reuses
# the workbook's existing date style rather than creating a brand-new one.
$birds.Range("G12").Copy()
$birds.Range("G13").PasteSpecial(-4122)
$birds.Cells.Item(13, 7).Value = $recordDate
$birds.Cells.Item(13, 8).Value = "a3"
$birds.Cells.Item(13, 9).Value = 777

# Row 14: same species/subspecies/gender, with mother (E) / father (F) ids, cage a3, user 777
$birds.Cells.Item(14, 1).Value = 13
$birds.Cells.Item(14, 2).Value = "American Gouldian"
$birds.Cells.Item(14, 3).Value = "Central America"
$birds.Cells.Item(14, 4).Value = "Male"
$birds.Cells.Item(14, 5).Value = 7
$birds.Cells.Item(14, 6).Value = 12
$birds.Range("G12").Copy()
$birds.Range("G14").PasteSpecial(-4122)
$birds.Cells.Item(14, 7).Value = $recordDate
$birds.Cells.Item(14, 8).Value = "a3"
$birds.Cells.Item(14, 9).Value = 777

$wb.Save()
